# Updates the crypto price/volume table to the latest scraped snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces Excel to store the Price column as text (it
# already is, in the source data) instead of auto-converting numeric-
# looking strings like "5.40" into the number 5.4. Resetting the style to
# "Normal" afterwards clears the quote-prefix formatting so cell styling
# is left exactly as it was.
function Set-TextValue($range, $value) {
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '38.326.50'
$ws.Range("E2").Value = '  +1.67%  '
Set-TextValue $ws.Range("D3") '2.078.11'
$ws.Range("E3").Value = '  +2.30%  '
$ws.Range("E4").Value = '  +0.12%  '
Set-TextValue $ws.Range("D5") '227.55'
$ws.Range("E5").Value = '  -0.05%  '
Set-TextValue $ws.Range("D6") '0.611'
$ws.Range("E6").Value = '  +0.87%  '
Set-TextValue $ws.Range("D7") '60.44'
$ws.Range("E7").Value = '  +0.38%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("E9").Value = '  +1.84%  '
Set-TextValue $ws.Range("D10") '0.0832'
$ws.Range("E10").Value = '  +1.13%  '
$ws.Range("E11").Value = '  -0.01%  '
Set-TextValue $ws.Range("D12") '2.390.13'
$ws.Range("E12").Value = '  +2.17%  '
Set-TextValue $ws.Range("D13") '14.77'
$ws.Range("E13").Value = '  +1.93%  '
Set-TextValue $ws.Range("D14") '22.34'
$ws.Range("E14").Value = '  +6.80%  '
Set-TextValue $ws.Range("D15") '0.782'
$ws.Range("E15").Value = '  +1.43%  '
Set-TextValue $ws.Range("D16") '5.40'
$ws.Range("E16").Value = '  +3.59%  '
Set-TextValue $ws.Range("D17") '2.073.70'
$ws.Range("E17").Value = '  +1.92%  '
Set-TextValue $ws.Range("D18") '38.258.68'
$ws.Range("E18").Value = '  +1.75%  '
Set-TextValue $ws.Range("D19") '71.41'
$ws.Range("E19").Value = '  +3.00%  '
Set-TextValue $ws.Range("D20") '6.01'
$ws.Range("E20").Value = '  +1.80%  '
Set-TextValue $ws.Range("D21") '0.0₃0831'
$ws.Range("E21").Value = '  +1.42%  '
Set-TextValue $ws.Range("D22") '225.05'
$ws.Range("E22").Value = '  +0.52%  '
$ws.Range("E23").Value = '  -0.03%  '
Set-TextValue $ws.Range("D24") '2.43'
$ws.Range("E24").Value = '  +0.39%  '
$ws.Range("E25").Value = '  +1.57%  '
Set-TextValue $ws.Range("D26") '169.54'
$ws.Range("E26").Value = '  +0.87%  '
Set-TextValue $ws.Range("D27") '9.39'
$ws.Range("E27").Value = '  +0.97%  '
Set-TextValue $ws.Range("D28") '0.135'
$ws.Range("E28").Value = '  +4.61%  '
Set-TextValue $ws.Range("D29") '18.99'
$ws.Range("E29").Value = '  +1.64%  '
$ws.Range("E30").Value = '  +8.63%  '
$ws.Range("E31").Value = '  -0.15%  '
Set-TextValue $ws.Range("D32") '2.33'
$ws.Range("E32").Value = '  +5.87%  '
$ws.Range("E33").Value = '  +7.45%  '
Set-TextValue $ws.Range("D34") '4.48'
$ws.Range("E34").Value = '  +2.39%  '
Set-TextValue $ws.Range("D35") '0.0604'
$ws.Range("E35").Value = '  +0.47%  '
# Row 36 becomes LidoDAOToken (previously THORChain); THORChain moves to row 37.
$ws.Range("B36").Value = 'LidoDAOToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Range("D36") '2.37'
$ws.Range("E36").Value = '  +2.20%  '
# Row 37 becomes THORChain (previously LidoDAOToken).
$ws.Range("B37").Value = 'THORChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue $ws.Range("D37") '6.38'
$ws.Range("E37").Value = '  -2.04%  '

Set-TextValue $ws.Range("D38") '3.53'
$ws.Range("E38").Value = '  +3.30%  '
$ws.Range("E39").Value = '  -0.07%  '
Set-TextValue $ws.Range("D40") '18.29'
$ws.Range("E40").Value = '  +1.40%  '
Set-TextValue $ws.Range("D41") '1.535.06'
$ws.Range("E41").Value = '  +0.32%  '
Set-TextValue $ws.Range("D42") '99.78'
$ws.Range("E42").Value = '  +3.58%  '
$ws.Range("E43").Value = '  +2.17%  '
Set-TextValue $ws.Range("D44") '0.0922'
$ws.Range("E44").Value = '  +1.66%  '
$ws.Range("E45").Value = '  -1.43%  '
Set-TextValue $ws.Range("D46") '7.70'
$ws.Range("E46").Value = '  +9.65%  '
Set-TextValue $ws.Range("D47") '4.10'
$ws.Range("E47").Value = '  -0.20%  '
$ws.Range("E48").Value = '  +0.88%  '
$ws.Range("E49").Value = '  +2.47%  '
Set-TextValue $ws.Range("D50") '2.97'
$ws.Range("E50").Value = '  +1.40%  '
Set-TextValue $ws.Range("D51") '2.276.71'
$ws.Range("E51").Value = '  +2.30%  '
